$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme
Write-Output ($cs | Get-Member -Force | Out-String)
